$wb = $excel.ActiveWorkbook

# Update the Aggr_generation sheet with new figures for the 2030 base case
$ws = $wb.Worksheets.Item("Aggr_generation")

$ws.Range("B2").Value = 8219
$ws.Range("B3").Value = 1496
$ws.Range("B4").Value = 5106
$ws.Range("B5").Value = 2412
$ws.Range("B6").Value = 4302
$ws.Range("B7").Value = 2627
$ws.Range("B8").Value = 3344
$ws.Range("B9").Value = 6646
$ws.Range("B10").Value = 9863
$ws.Range("B11").Value = 2219

# Make "Aggr_generation" the active sheet/tab
$ws.Activate()
